$d = $word.ActiveDocument

$pairs = @(
    @("721÷6=", "553÷2="),
    @("545÷6=", "466÷3="),
    @("719÷9=", "682÷3="),
    @("911÷5=", "715÷4="),
    @("616÷5=", "834÷6="),
    @("951÷4=", "763÷2="),
    @("741÷7=", "577÷7="),
    @("896÷6=", "889÷4="),
    @("218÷3=", "774÷3="),
    @("875÷7=", "871÷6="),
    @("800÷3=", "764÷7="),
    @("937÷3=", "843÷5="),
    @("706÷4=", "309÷3="),
    @("917÷8=", "701÷9="),
    @("450÷4=", "122÷5="),
    @("868÷4=", "651÷5="),
    @("285÷6=", "779÷5="),
    @("697÷9=", "261÷6="),
    @("308÷5=", "553÷6="),
    @("253÷9=", "372÷9="),
    @("892÷6=", "265÷2="),
    @("943÷5=", "689÷8="),
    @("920÷5=", "158÷4="),
    @("513÷7=", "962÷3="),
    @("715÷3=", "303÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
